$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the optical-flow method value from "mog2" to "gmg" (B2:C2)
$ws.Range("B2").Value = "gmg"
$ws.Range("C2").Value = "gmg"

# Add a new parameter row for the GMG comparison threshold
$ws.Range("A29").Value = "CompResThres"
$ws.Range("B29").Value = 0.05
$ws.Range("C29").Value = 0.05

# Document the new parameter with a cell comment
$ws.Range("A29").AddComment("Default 0.05`nGrimson-Stauffer 0.0") | Out-Null

# Column A needs to be a bit wider to fit the new label; the rest keep
# their previous (default) width
$ws.Columns.Item(1).ColumnWidth = 13.8877551020408

# Move the active selection back up to C2
$ws.Range("C2").Select() | Out-Null
